$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header in H1 - copy the formatting of the existing header row (G1)
# so the new column matches the bold/centered/bordered header style, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Populate the new "Save" column values for each data row (2-9)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("H9").Value = 0
